$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.216.02"
$ws.Range("E2").Value = "  +0.68%  "
$ws.Range("D3").Value = "1.852.33"
$ws.Range("E3").Value = "  +1.22%  "
$ws.Range("E4").Value = "  -0.39%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.68"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.63%  "
$ws.Range("E6").Value = "  -0.37%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4635"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3711"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07284"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.87%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8864"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.89%  "
$ws.Range("E11").Value = "  +1.51%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07866"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.27%  "
$ws.Range("D13").Value = "1.815.43"
$ws.Range("E13").Value = "  +2.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.389"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.90%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.510"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.44%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.01"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.56%  "
$ws.Range("E17").Value = "  -0.43%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008931"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.81%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.001"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.67"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.80%  "
$ws.Range("D21").Value = "27.240.04"
$ws.Range("E21").Value = "  +0.70%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.080"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.57%  "
$ws.Range("E23").Value = "  -0.21%  "
$ws.Range("D24").Value = "2.072.89"
$ws.Range("E24").Value = "  +3.62%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.951"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +5.67%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.23"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.06%  "
$ws.Range("E27").Value = "  -0.59%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.045"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.06%  "
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.039"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.72%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08801"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.02%  "
$ws.Range("E32").Value = "  +6.12%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7688"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +5.60%  "
$ws.Range("E34").Value = "  +2.80%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.521"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +1.76%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.720"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +10.38%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.111"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +3.15%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01938"
$ws.Range("D38").ClearFormats()
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05217"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.33%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.936"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.58%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.041"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.99%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5123"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1628"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.12%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.446"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +3.23%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4797"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.94%  "
$ws.Range("E46").Value = "  +2.21%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.000"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.43%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "102.89"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.49%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.643"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.73%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06208"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "65.55"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.26%  "
